$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '26.024.24'
$ws.Range('E2').Value = '  +0.40%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '1.643.66'
$ws.Range('E3').Value = '  +0.40%  '
$ws.Range('E4').Value = '  +0.53%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '215.95'
$ws.Range('E5').Value = '  +0.65%  '
$ws.Range('E6').Value = '  +0.19%  '
$ws.Range('E8').Value = '  +0.43%  '
$ws.Range('E9').Value = '  +0.57%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '19.59'
$ws.Range('E10').Value = '  +0.03%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.0797'
$ws.Range('E11').Value = '  +0.37%  '
$ws.Range('E12').Value = '  +0.63%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '1.655.18'
$ws.Range('E13').Value = '  +0.63%  '
$ws.Range('E14').Value = '  +0.29%  '
$ws.Range('B15').Value = 'Litecoin'
$ws.Range('C15').Value = 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc'
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '63.42'
$ws.Range('E15').Value = '  +1.41%  '
$ws.Range('B16').Value = 'ShibaInu'
$ws.Range('C16').Value = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '0.0₃0762'
$ws.Range('E16').Value = '  +0.96%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '26.051.99'
$ws.Range('E17').Value = '  +0.44%  '
$ws.Range('E18').Value = '  +0.46%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '194.31'
$ws.Range('E19').Value = '  +0.33%  '
$ws.Range('E20').Value = '  -0.71%  '
$ws.Range('E21').Value = '  +0.16%  '
$ws.Range('E22').Value = '  -1.21%  '
$ws.Range('E23').Value = '  +4.42%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '1.80'
$ws.Range('E24').Value = '  -1.09%  '
$ws.Range('E25').Value = '  +0.43%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '143.20'
$ws.Range('E26').Value = '  -0.37%  '
$ws.Range('E27').Value = '  +0.57%  '
$ws.Range('E28').Value = '  +0.58%  '
$ws.Range('E29').Value = '  +0.52%  '
$ws.Range('E30').Value = '  -1.09%  '
$ws.Range('E31').Value = '  -0.18%  '
$ws.Range('E32').Value = '  +1.41%  '
$ws.Range('E33').Value = '  -0.81%  '
$ws.Range('E34').Value = '  +1.59%  '
$ws.Range('E35').Value = '  +0.34%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '1.131.12'
$ws.Range('E36').Value = '  -0.60%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '0.541'
$ws.Range('E37').Value = '  -0.71%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '2.47'
$ws.Range('E38').Value = '  +0.39%  '
$ws.Range('E39').Value = '  +0.20%  '
$ws.Range('E40').Value = '  +0.73%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '99.14'
$ws.Range('E41').Value = '  -0.28%  '
$ws.Range('E42').Value = '  -0.02%  '
$ws.Range('E43').Value = '  +2.20%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '56.43'
$ws.Range('E44').Value = '  +0.00%  '
$ws.Range('E46').Value = '  +2.91%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '7.75'
$ws.Range('E47').Value = '  +1.26%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '0.415'
$ws.Range('E48').Value = '  -0.13%  '
$ws.Range('E49').Value = '  +0.25%  '
$ws.Range('E50').Value = '  -0.78%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '1.18'
$ws.Range('E51').Value = '  +3.07%  '
